$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the solution descriptions in column A (row 3-9) with the new
# human-readable Spanish text instead of the raw boolean-list repr.
$ws.Range("A3").Value = "Se almacenan los objetos: 1, 3"
$ws.Range("A4").Value = "Se almacenan los objetos: 1, 2"
$ws.Range("A5").Value = "Se almacenan los objetos: 2, 3"
$ws.Range("A6").Value = "Se almacena el objeto: 1"
$ws.Range("A7").Value = "Se almacena el objeto: 3"
$ws.Range("A8").Value = "Se almacena el objeto: 2"
$ws.Range("A9").Value = "No se almacenan objetos"

# The text cells in column A keep their thin border but drop the
# center alignment they used to share with column B (a new style,
# border-only, is introduced for them), while column B keeps its
# original bordered + centered look for the numeric totals.
$rngA = $ws.Range("A3:A9")
$rngA.ClearFormats()
$rngA.Borders.ColorIndex = 1
